# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The detail table (rows 16-42, columns B:G) is rewritten: the data that used
# to be grouped by worker (each worker's 9 "Periodo Mora" rows together, in
# descending period order) is now grouped by period (ascending order), with
# the three workers repeated inside every period group. The "Valor Mora" of
# 28090 stays attached to period 2105 (now the last group instead of the
# first), every other period keeps 35112, and "Salario Basico" (G) is
# untouched at 877803 throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=16; B="CC"; C="8853553";         D="EDINSON SARMIENTO MUÑOZ";        E="2009"; F=35112; G=877803 },
    @{ Row=17; B="CC"; C="79655757";        D="HENRY ALEXANDER HERNANDEZ CUERVO"; E="2009"; F=35112; G=877803 },
    @{ Row=18; B="PE"; C="802151520031995"; D="ANDERSON JOSE FERNANDEZ NIERES";  E="2009"; F=35112; G=877803 },

    @{ Row=19; B="CC"; C="8853553";         D="EDINSON SARMIENTO MUÑOZ";        E="2010"; F=35112; G=877803 },
    @{ Row=20; B="CC"; C="79655757";        D="HENRY ALEXANDER HERNANDEZ CUERVO"; E="2010"; F=35112; G=877803 },
    @{ Row=21; B="PE"; C="802151520031995"; D="ANDERSON JOSE FERNANDEZ NIERES";  E="2010"; F=35112; G=877803 },

    @{ Row=22; B="CC"; C="8853553";         D="EDINSON SARMIENTO MUÑOZ";        E="2011"; F=35112; G=877803 },
    @{ Row=23; B="CC"; C="79655757";        D="HENRY ALEXANDER HERNANDEZ CUERVO"; E="2011"; F=35112; G=877803 },
    @{ Row=24; B="PE"; C="802151520031995"; D="ANDERSON JOSE FERNANDEZ NIERES";  E="2011"; F=35112; G=877803 },

    @{ Row=25; B="CC"; C="8853553";         D="EDINSON SARMIENTO MUÑOZ";        E="2012"; F=35112; G=877803 },
    @{ Row=26; B="CC"; C="79655757";        D="HENRY ALEXANDER HERNANDEZ CUERVO"; E="2012"; F=35112; G=877803 },
    @{ Row=27; B="PE"; C="802151520031995"; D="ANDERSON JOSE FERNANDEZ NIERES";  E="2012"; F=35112; G=877803 },

    @{ Row=28; B="CC"; C="8853553";         D="EDINSON SARMIENTO MUÑOZ";        E="2101"; F=35112; G=877803 },
    @{ Row=29; B="CC"; C="79655757";        D="HENRY ALEXANDER HERNANDEZ CUERVO"; E="2101"; F=35112; G=877803 },
    @{ Row=30; B="PE"; C="802151520031995"; D="ANDERSON JOSE FERNANDEZ NIERES";  E="2101"; F=35112; G=877803 },

    @{ Row=31; B="CC"; C="8853553";         D="EDINSON SARMIENTO MUÑOZ";        E="2102"; F=35112; G=877803 },
    @{ Row=32; B="CC"; C="79655757";        D="HENRY ALEXANDER HERNANDEZ CUERVO"; E="2102"; F=35112; G=877803 },
    @{ Row=33; B="PE"; C="802151520031995"; D="ANDERSON JOSE FERNANDEZ NIERES";  E="2102"; F=35112; G=877803 },

    @{ Row=34; B="CC"; C="8853553";         D="EDINSON SARMIENTO MUÑOZ";        E="2103"; F=35112; G=877803 },
    @{ Row=35; B="CC"; C="79655757";        D="HENRY ALEXANDER HERNANDEZ CUERVO"; E="2103"; F=35112; G=877803 },
    @{ Row=36; B="PE"; C="802151520031995"; D="ANDERSON JOSE FERNANDEZ NIERES";  E="2103"; F=35112; G=877803 },

    @{ Row=37; B="CC"; C="8853553";         D="EDINSON SARMIENTO MUÑOZ";        E="2104"; F=35112; G=877803 },
    @{ Row=38; B="CC"; C="79655757";        D="HENRY ALEXANDER HERNANDEZ CUERVO"; E="2104"; F=35112; G=877803 },
    @{ Row=39; B="PE"; C="802151520031995"; D="ANDERSON JOSE FERNANDEZ NIERES";  E="2104"; F=35112; G=877803 },

    @{ Row=40; B="CC"; C="8853553";         D="EDINSON SARMIENTO MUÑOZ";        E="2105"; F=28090; G=877803 },
    @{ Row=41; B="CC"; C="79655757";        D="HENRY ALEXANDER HERNANDEZ CUERVO"; E="2105"; F=28090; G=877803 },
    @{ Row=42; B="PE"; C="802151520031995"; D="ANDERSON JOSE FERNANDEZ NIERES";  E="2105"; F=28090; G=877803 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
